# "Generate Report for Handoff" - refresh the report's source-file GUID,
# localized-artifact hash, and handoff timestamps for the new run.

$wb = $excel.ActiveWorkbook

$newGuid = "2816ee50-5f2f-488d-8155-d20c5f3c4e61"
$newHash = "b16b81ddd90bf4a143a7e996207fa76ddd38b7e3"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newOverviewHandoffDate = "2016-03-24 07:07:45"
$newZhHandoffDatetime   = "2016-03-24 07:07:40"

# The hyperlink targets themselves (the rels) keep pointing at the
# originally-generated artifacts; only the visible file names change.
$mdTarget    = "https://github.com/OpenLocalizationTest/oltest/blob/093716a7cefd2aee3dffb0091f2620760d613b0d/e2e/fcf1c24b-55a5-4919-80de-5c3878c14a7b.md"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/843f7badc0291e6b16acbb749f3e8f268d23a87a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fcf1c24b-55a5-4919-80de-5c3878c14a7b.92576864edfdc23edfc86040d33ad9fd9023f8b0.zh-cn.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6dab8a52b3dfed150f853f33cee5ad6b3b4f5ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fcf1c24b-55a5-4919-80de-5c3878c14a7b.92576864edfdc23edfc86040d33ad9fd9023f8b0.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMdName)
$wsOverview.Range("D2").Value = $newOverviewHandoffDate

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdTarget, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhXlfTarget, "", "", $newZhXlfName)
$wsZhCn.Range("E2").Value = $newZhHandoffDatetime

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdTarget, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deXlfTarget, "", "", $newDeXlfName)
